$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 17546834
$ws.Range("I76").Value = 20002944
$ws.Range("J76").Value = 3178.5715
$ws.Range("K76").Value = 20002944
$ws.Range("L76").Value = 3178.5715
$ws.Range("M76").Value = -20002629
$ws.Range("N76").Value = -3808.5715
$ws.Range("H79").Value = 17546834
$ws.Range("I79").Value = 20002944
$ws.Range("J79").Value = 3178.5715
$ws.Range("K79").Value = 20002944
$ws.Range("L79").Value = 3178.5715
$ws.Range("M79").Value = -20001852
$ws.Range("N79").Value = -5362.5715

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 23877.334
$ws.Range("I21").Value = 807.5
$ws.Range("K21").Value = 807.5
$ws.Range("M21").Value = -433.5
$ws.Range("H32").Value = 1280906.8
$ws.Range("I32").Value = 1493479
$ws.Range("J32").Value = 5473
$ws.Range("K32").Value = 1493479
$ws.Range("L32").Value = 5473
$ws.Range("M32").Value = -1493192
$ws.Range("N32").Value = -6047
$ws.Range("H102").Value = 2518.6667
$ws.Range("I102").Value = 2778
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2778
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -1156
$ws.Range("N102").Value = -5244
$ws.Range("H110").Value = 1045.4546
$ws.Range("I110").Value = 1178.4445
$ws.Range("J110").Value = 447
$ws.Range("K110").Value = 1178.4445
$ws.Range("L110").Value = 447
$ws.Range("M110").Value = 866.5554999999999
$ws.Range("N110").Value = -4537

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 55613516
$ws.Range("I20").Value = 33347344
$ws.Range("J20").Value = 83446230
$ws.Range("K20").Value = 33347344
$ws.Range("L20").Value = 83446230
$ws.Range("M20").Value = -33347097
$ws.Range("N20").Value = -83446724
$ws.Range("H86").Value = 1943.97
$ws.Range("I86").Value = 1951
$ws.Range("J86").Value = 1599.5
$ws.Range("K86").Value = 1951
$ws.Range("L86").Value = 1599.5
$ws.Range("M86").Value = -828
$ws.Range("N86").Value = -3845.5
$ws.Range("H89").Value = 1943.97
$ws.Range("I89").Value = 1951
$ws.Range("J89").Value = 1599.5
$ws.Range("K89").Value = 9755
$ws.Range("L89").Value = 7997.5
$ws.Range("M89").Value = -4139
$ws.Range("N89").Value = -19229.5
$ws.Range("H105").Value = 1554.6562
$ws.Range("I105").Value = 1623.8096
$ws.Range("J105").Value = 1422.6364
$ws.Range("K105").Value = 1623.8096
$ws.Range("L105").Value = 1422.6364
$ws.Range("M105").Value = 123.1904
$ws.Range("N105").Value = -4916.6364
$ws.Range("H107").Value = 2133.2222
$ws.Range("I107").Value = 2274.875
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2274.875
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -354.875
$ws.Range("N107").Value = -4840

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14927049
$ws.Range("I31").Value = 34484092
$ws.Range("J31").Value = 1937.7368
$ws.Range("K31").Value = 34484092
$ws.Range("L31").Value = 1937.7368
$ws.Range("M31").Value = -34483797
$ws.Range("N31").Value = -2527.7368
$ws.Range("H34").Value = 14927049
$ws.Range("I34").Value = 34484092
$ws.Range("J34").Value = 1937.7368
$ws.Range("K34").Value = 34484092
$ws.Range("L34").Value = 1937.7368
$ws.Range("M34").Value = -34483890
$ws.Range("N34").Value = -2341.7368

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2660179.8
$ws.Range("I5").Value = 2262959.8
$ws.Range("J5").Value = 3335453.8
$ws.Range("K5").Value = 6788879.399999999
$ws.Range("L5").Value = 10006361.4
$ws.Range("M5").Value = -6788767.399999999
$ws.Range("N5").Value = -10006585.4
$ws.Range("H68").Value = 3801.9736
$ws.Range("J68").Value = 8613.333000000001
$ws.Range("L68").Value = 25839.999
$ws.Range("N68").Value = -27461.999
$ws.Range("H71").Value = 3801.9736
$ws.Range("J71").Value = 8613.333000000001
$ws.Range("L71").Value = 77519.997
$ws.Range("N71").Value = -85631.997
$ws.Range("H131").Value = 26707.025
$ws.Range("J131").Value = 1389
$ws.Range("L131").Value = 4167
$ws.Range("N131").Value = -14247
$ws.Range("H132").Value = 1569.6177
$ws.Range("I132").Value = 1189
$ws.Range("J132").Value = 1777.2273
$ws.Range("K132").Value = 10701
$ws.Range("L132").Value = 15995.0457
$ws.Range("M132").Value = -8171
$ws.Range("N132").Value = -21055.0457
$ws.Range("H135").Value = 2660179.8
$ws.Range("I135").Value = 2262959.8
$ws.Range("J135").Value = 3335453.8
$ws.Range("K135").Value = 20366638.2
$ws.Range("L135").Value = 30019084.2
$ws.Range("M135").Value = -20364103.2
$ws.Range("N135").Value = -30024154.2
$ws.Range("H136").Value = 1860.7858
$ws.Range("I136").Value = 1285.1
$ws.Range("J136").Value = 3300
$ws.Range("K136").Value = 3855.3
$ws.Range("L136").Value = 9900
$ws.Range("M136").Value = 1244.7
$ws.Range("N136").Value = -20100
$ws.Range("H137").Value = 3063.577
$ws.Range("I137").Value = 1881.3334
$ws.Range("J137").Value = 4675.727
$ws.Range("K137").Value = 5644.0002
$ws.Range("L137").Value = 14027.181
$ws.Range("M137").Value = -544.0002000000004
$ws.Range("N137").Value = -24227.181
$ws.Range("H138").Value = 1686.1852
$ws.Range("I138").Value = 781.9524
$ws.Range("K138").Value = 2345.8572
$ws.Range("M138").Value = 2794.1428
$ws.Range("H139").Value = 42955.5
$ws.Range("I139").Value = 50540
$ws.Range("K139").Value = 151620
$ws.Range("M139").Value = -146480
$ws.Range("H140").Value = 3229.7083
$ws.Range("I140").Value = 3885
$ws.Range("J140").Value = 2312.3
$ws.Range("K140").Value = 11655
$ws.Range("L140").Value = 6936.900000000001
$ws.Range("M140").Value = -6475
$ws.Range("N140").Value = -17296.9
$ws.Range("H141").Value = 2026.7693
$ws.Range("I141").Value = 1862.3334
$ws.Range("K141").Value = 5587.0002
$ws.Range("M141").Value = -407.0002000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = $null
$ws.Range("H27").Value = 1000
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1332
$ws.Range("H70").Value = 4313761.5
$ws.Range("I70").Value = 2141688.5
$ws.Range("K70").Value = 2141688.5
$ws.Range("M70").Value = -2141418.5
$ws.Range("H73").Value = 4313761.5
$ws.Range("I73").Value = 2141688.5
$ws.Range("K73").Value = 2141688.5
$ws.Range("M73").Value = -2140752.5
$ws.Range("H132").Value = 10786637
$ws.Range("I132").Value = 6970459.5
$ws.Range("J132").Value = 23380022
$ws.Range("K132").Value = 20911378.5
$ws.Range("L132").Value = 70140066
$ws.Range("M132").Value = -20908848.5
$ws.Range("N132").Value = -70145126

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = $null
$ws.Range("H40").Value = 3475687.5
$ws.Range("J40").Value = 4508.3
$ws.Range("L40").Value = 4508.3
$ws.Range("N40").Value = -4780.3

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 3513.3
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 3513.3
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 3513.3
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = -3993.3
$ws.Range("H41").Value = 4000
$ws.Range("J41").Value = 4000
$ws.Range("L41").Value = 4000
$ws.Range("N41").Value = -4780
$ws.Range("H107").Value = 23411.555
$ws.Range("I107").Value = 25975.5
$ws.Range("K107").Value = 77926.5
$ws.Range("M107").Value = -76006.5
$ws.Range("H126").Value = 21826128
$ws.Range("I126").Value = 24554270
$ws.Range("K126").Value = 73662810
$ws.Range("M126").Value = -73660340
